$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 63, shifting the existing
# rows 63-94 down to 65-96 (formatting is inherited from the row above,
# which keeps column D's date number-format on the new rows too).
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()

# New row 63: Lapins / Primera, Región de Ñuble
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = "Vega Monumental Concepción"
$ws.Range("C63").Value = "Bíobío"
$ws.Range("D63").Value = 44572
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100103
$ws.Range("H63").Value = "Frutos de hueso (carozo)"
$ws.Range("I63").Value = 100103001
$ws.Range("J63").Value = "Cereza"
$ws.Range("K63").Value = "Lapins"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 200
$ws.Range("N63").Value = 5500
$ws.Range("O63").Value = 6000
$ws.Range("P63").Value = 5750
$ws.Range("Q63").Value = "$/bandeja 10 kilos"
$ws.Range("R63").Value = "Región de Ñuble"
$ws.Range("S63").Value = 575
$ws.Range("T63").Value = 10

# New row 64: Lapins / Segunda, Región de Ñuble
$ws.Range("A64").Value = 11
$ws.Range("B64").Value = "Vega Monumental Concepción"
$ws.Range("C64").Value = "Bíobío"
$ws.Range("D64").Value = 44572
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100103
$ws.Range("H64").Value = "Frutos de hueso (carozo)"
$ws.Range("I64").Value = 100103001
$ws.Range("J64").Value = "Cereza"
$ws.Range("K64").Value = "Lapins"
$ws.Range("L64").Value = "Segunda"
$ws.Range("M64").Value = 100
$ws.Range("N64").Value = 5000
$ws.Range("O64").Value = 5000
$ws.Range("P64").Value = 5000
$ws.Range("Q64").Value = "$/bandeja 10 kilos"
$ws.Range("R64").Value = "Región de Ñuble"
$ws.Range("S64").Value = 500
$ws.Range("T64").Value = 10
